# Commit: "Added Pause Screen to MainGame" -
# "added a pause screen and score panel to the MainGame scene."
#
# In this workbook that translates to a single content edit: the feature
# description in the task list was reworded from
#   "The ability for the charcter to move in cardinal directions and Jump
#    (Camera is vertically fixed on player"
# to
#   "The ability for the charcter to move/jump in cardinal directions
#    (Camera is vertically fixed on player"
# (cell A10 on the "Project Plan" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan")

$ws.Range("A10").Value = "The ability for the charcter to move/jump in cardinal directions (Camera is vertically fixed on player"

# Match the author's final selection/view state in the sheet.
$ws.Range("A10").Select()
